$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the width/height values (columns B and C) for rows 2-6 back to 1
$ws.Range("B2:C6").Value = 1

# Update the active selection to reflect where the user last clicked
$ws.Range("B6").Select()
